# Applies the diff: adds two new sheets ("o_20", "o_20_jumbled") with
# the same header row as "o_10", adds a new "evaluator_partial_correctness"
# column to "o_10", updates the llm_response text, and adds a partial
# correctness output value.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("o_10")

# --- Update existing sheet "o_10" ---

# New header cell E1: copy the formatting of the other header cells
# (bold font, thin border, centered/top aligned) then set its text.
$ws1.Range("A1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("E1").Value = "evaluator_partial_correctness"

# Update C2 (llm_response) text to the new multi-line value
$ws1.Range("C2").Value = "The shortest path from node A to node J is:`nA -> B -> C -> D -> E -> F -> G -> H -> I -> J"
# Undo the row auto-height bump the embedded line break triggers, so row 2
# keeps relying on the sheet's default row height like before.
$ws1.Rows(2).EntireRow.AutoFit()

# New E2 cell value (evaluator_partial_correctness)
$ws1.Range("E2").Value = "Output: 10/10"

$headers = @("prompt", "solution", "llm_response", "evaluator_response", "evaluator_partial_correctness")

# --- Add new sheets (same header row as "o_10"), in order right after it ---

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "o_20"
$ws1.Range("A1:E1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "o_20_jumbled"
$ws1.Range("A1:E1").Copy()
$ws3.Range("A1:E1").PasteSpecial(-4122)  # xlPasteFormats
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers[$i]
}

$excel.CutCopyMode = $false
$ws1.Activate()
[void]$ws1.Range("A1").Select()
